# Append a new data row (97) after the existing last row (96) on the
# active worksheet, matching the source/author/core inflation figures
# for the added period.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A97").Value = 1.776
$ws.Range("B97").Value = 1.009
$ws.Range("C97").Value = 1.644
